# Atualização de bases das ligas, do dia: 28-05-2024 às 07:50
#
# The source feed re-paired several fixture rows (e.g. two matches sharing
# the same Date got re-ordered by a refreshed match id), which manifests in
# the sheet as pairs of adjacent rows whose B:AD contents (match id through
# the closing-line odds columns) are fully swapped while the leading "id"
# column (A) stays anchored to its row.
#
# Apply that swap for every affected row pair.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Belgium First Division A")

# (row1, row2) pairs whose B:AD data got swapped
$pairs = @(
    @(159, 160),
    @(164, 165),
    @(175, 176),
    @(181, 182),
    @(187, 188),
    @(190, 191),
    @(313, 314)
)

foreach ($pair in $pairs) {
    $r1 = $pair[0]
    $r2 = $pair[1]

    $rng1 = $ws.Range("B$r1`:AD$r1")
    $rng2 = $ws.Range("B$r2`:AD$r2")

    $vals1 = $rng1.Value2
    $vals2 = $rng2.Value2

    $rng1.Value = $vals2
    $rng2.Value = $vals1
}
